$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1723.0435
$ws.Range("I98").Value = 1231.5
$ws.Range("J98").Value = 5000
$ws.Range("K98").Value = 1231.5
$ws.Range("L98").Value = 5000
$ws.Range("M98").Value = 266.5
$ws.Range("N98").Value = -7996

$ws.Range("H122").Value = 1723.0435
$ws.Range("I122").Value = 1231.5
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 3694.5
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -1244.5
$ws.Range("N122").Value = -19900

$ws.Range("H132").Value = 1235.4166
$ws.Range("I132").Value = 1302.6786
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 3908.0358
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = -1378.0358
$ws.Range("N132").Value = -8060

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 759.0862
$ws.Range("I2").Value = 684.5102000000001
$ws.Range("J2").Value = 1165.1111
$ws.Range("K2").Value = 684.5102000000001
$ws.Range("L2").Value = 1165.1111
$ws.Range("M2").Value = -571.5102000000001
$ws.Range("N2").Value = -1391.1111

$ws.Range("H32").Value = 7681.255
$ws.Range("I32").Value = 3698.8416
$ws.Range("J32").Value = 28091.125
$ws.Range("K32").Value = 3698.8416
$ws.Range("L32").Value = 28091.125
$ws.Range("M32").Value = -3411.8416
$ws.Range("N32").Value = -28665.125

$ws.Range("H45").Value = 1127.5625
$ws.Range("I45").Value = 1031
$ws.Range("J45").Value = 1340
$ws.Range("K45").Value = 1031
$ws.Range("L45").Value = 1340
$ws.Range("M45").Value = -654
$ws.Range("N45").Value = -2094

$ws.Range("H74").Value = 957.7377300000001
$ws.Range("I74").Value = 736.93475
$ws.Range("J74").Value = 1634.8667
$ws.Range("K74").Value = 736.93475
$ws.Range("L74").Value = 1634.8667
$ws.Range("M74").Value = 137.06525
$ws.Range("N74").Value = -3382.8667

$ws.Range("H77").Value = 957.7377300000001
$ws.Range("I77").Value = 736.93475
$ws.Range("J77").Value = 1634.8667
$ws.Range("K77").Value = 3684.67375
$ws.Range("L77").Value = 8174.333500000001
$ws.Range("M77").Value = 683.3262500000001
$ws.Range("N77").Value = -16910.3335

$ws.Range("H116").Value = 759.0862
$ws.Range("I116").Value = 684.5102000000001
$ws.Range("J116").Value = 1165.1111
$ws.Range("K116").Value = 684.5102000000001
$ws.Range("L116").Value = 1165.1111
$ws.Range("M116").Value = 1609.4898
$ws.Range("N116").Value = -5753.1111

$ws.Range("H122").Value = 1444.0488
$ws.Range("I122").Value = 1228
$ws.Range("J122").Value = 1966.1666
$ws.Range("K122").Value = 3684
$ws.Range("L122").Value = 5898.4998
$ws.Range("M122").Value = -1234
$ws.Range("N122").Value = -10798.4998

$ws.Range("H132").Value = 2430.2693
$ws.Range("I132").Value = 1844.8206
$ws.Range("J132").Value = 4186.615
$ws.Range("K132").Value = 5534.4618
$ws.Range("L132").Value = 12559.845
$ws.Range("M132").Value = -3004.4618
$ws.Range("N132").Value = -17619.845

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 759.0862
$ws.Range("I3").Value = 684.5102000000001
$ws.Range("J3").Value = 1165.1111
$ws.Range("K3").Value = 684.5102000000001
$ws.Range("L3").Value = 1165.1111
$ws.Range("M3").Value = -570.5102000000001
$ws.Range("N3").Value = -1393.1111

$ws.Range("H99").Value = 834.40424
$ws.Range("I99").Value = 866.8276
$ws.Range("J99").Value = 782.1667
$ws.Range("K99").Value = 866.8276
$ws.Range("L99").Value = 782.1667
$ws.Range("M99").Value = 631.1724
$ws.Range("N99").Value = -3778.1667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2424.1143
$ws.Range("I122").Value = 2756.158
$ws.Range("J122").Value = 2029.8125
$ws.Range("K122").Value = 8268.474
$ws.Range("L122").Value = 6089.4375
$ws.Range("M122").Value = -5818.474
$ws.Range("N122").Value = -10989.4375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H48").Value = 37039704
$ws.Range("J48").Value = 37039704
$ws.Range("L48").Value = 111119112
$ws.Range("N48").Value = -111119612

$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").ClearContents()

$ws.Range("H56").Value = 5045.5557
$ws.Range("I56").Value = 5045.5557
$ws.Range("K56").Value = 5045.5557
$ws.Range("M56").Value = -4515.5557

$ws.Range("H57").Value = 4065.2
$ws.Range("I57").Value = 440
$ws.Range("J57").Value = 9503
$ws.Range("K57").Value = 1320
$ws.Range("L57").Value = 28509
$ws.Range("M57").Value = -761
$ws.Range("N57").Value = -29627

$ws.Range("H58").Value = 2198.3333
$ws.Range("I58").Value = 1641.0714
$ws.Range("J58").Value = 10000
$ws.Range("K58").Value = 4923.2142
$ws.Range("L58").Value = 30000
$ws.Range("M58").Value = -4795.2142
$ws.Range("N58").Value = -30256

$ws.Range("H59").Value = 12400
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 12400
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 37200
$ws.Range("M59").ClearContents()
$ws.Range("N59").Value = -38280

$ws.Range("H60").Value = 1837.3334
$ws.Range("I60").Value = 500
$ws.Range("K60").Value = 1500
$ws.Range("M60").Value = -1249

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2121.1843
$ws.Range("I122").Value = 939.7895
$ws.Range("J122").Value = 3302.5789
$ws.Range("K122").Value = 2819.3685
$ws.Range("L122").Value = 9907.736699999999
$ws.Range("M122").Value = -369.3685
$ws.Range("N122").Value = -14807.7367

$ws.Range("H126").Value = 1616.5454
$ws.Range("I126").Value = 1616.5454
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 4849.6362
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -2379.6362
$ws.Range("N126").ClearContents()

$ws.Range("H132").Value = 2624.484
$ws.Range("I132").Value = 2373.1667
$ws.Range("J132").Value = 3486.1428
$ws.Range("K132").Value = 7119.500100000001
$ws.Range("L132").Value = 10458.4284
$ws.Range("M132").Value = -4589.500100000001
$ws.Range("N132").Value = -15518.4284

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1842.3529
$ws.Range("I7").Value = 1471.4286
$ws.Range("J7").Value = 2102
$ws.Range("K7").Value = 1471.4286
$ws.Range("L7").Value = 2102
$ws.Range("M7").Value = -1359.4286
$ws.Range("N7").Value = -2326

$ws.Range("H40").Value = 2247.111
$ws.Range("I40").Value = 1829.8334
$ws.Range("K40").Value = 1829.8334
$ws.Range("M40").Value = -1693.8334

$ws.Range("H126").Value = 1842.3529
$ws.Range("I126").Value = 1471.4286
$ws.Range("J126").Value = 2102
$ws.Range("K126").Value = 4414.2858
$ws.Range("L126").Value = 6306
$ws.Range("M126").Value = -1944.2858
$ws.Range("N126").Value = -11246

$ws.Range("H132").Value = 3831.9265
$ws.Range("I132").Value = 2231.4822
$ws.Range("J132").Value = 11300.667
$ws.Range("K132").Value = 6694.446599999999
$ws.Range("L132").Value = 33902.001
$ws.Range("M132").Value = -4164.446599999999
$ws.Range("N132").Value = -38962.001

$ws.Range("H136").Value = 2803.535
$ws.Range("I136").Value = 2138.963
$ws.Range("J136").Value = 3925
$ws.Range("K136").Value = 6416.889000000001
$ws.Range("L136").Value = 11775
$ws.Range("M136").Value = -3866.889000000001
$ws.Range("N136").Value = -16875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 334465.9
$ws.Range("I122").Value = 500930.34
$ws.Range("J122").Value = 1537
$ws.Range("K122").Value = 1502791.02
$ws.Range("L122").Value = 4611
$ws.Range("M122").Value = -1500341.02
$ws.Range("N122").Value = -9511

$ws.Range("H126").Value = 418055.4
$ws.Range("I126").Value = 556741.9
$ws.Range("J126").Value = 1996
$ws.Range("K126").Value = 1670225.7
$ws.Range("L126").Value = 5988
$ws.Range("M126").Value = -1667755.7
$ws.Range("N126").Value = -10928

$ws.Range("H132").Value = 979.1667
$ws.Range("I132").Value = 677.5769
$ws.Range("J132").Value = 1469.25
$ws.Range("K132").Value = 2032.7307
$ws.Range("L132").Value = 4407.75
$ws.Range("M132").Value = 497.2692999999999
$ws.Range("N132").Value = -9467.75

$ws.Range("H136").Value = 618.2917
$ws.Range("I136").Value = 527.27905
$ws.Range("J136").Value = 1401
$ws.Range("K136").Value = 1581.83715
$ws.Range("L136").Value = 4203
$ws.Range("M136").Value = 968.1628500000002
$ws.Range("N136").Value = -9303
